$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Shift dates for a batch of existing rows by +30 days (rows 501 and 509 are left untouched) ---
$ws.Range("A491").Value = 44840
$ws.Range("A492").Value = 44840
$ws.Range("A493").Value = 44840
$ws.Range("A494").Value = 44840

$ws.Range("A495").Value = 44841
$ws.Range("A496").Value = 44841
$ws.Range("A497").Value = 44841
$ws.Range("A498").Value = 44841
$ws.Range("A499").Value = 44841
$ws.Range("A500").Value = 44841

$ws.Range("A502").Value = 44842
$ws.Range("A503").Value = 44842
$ws.Range("A504").Value = 44842
$ws.Range("A505").Value = 44842
$ws.Range("A506").Value = 44842
$ws.Range("A507").Value = 44842
$ws.Range("A508").Value = 44842
$ws.Range("A510").Value = 44842

# --- Append new rows 511-519 ---
$ws.Range("A511").Value = 44844
$ws.Range("B511").Value = "HR51AR4836"
$ws.Range("C511").Value = "I10"
$ws.Range("D511").Value = "TIE MEMBER                "
$ws.Range("E511").Value = "WORK DONE DELIVERED"
$ws.Range("F511").Value = 6205
$ws.Range("G511").Value = "P PAY"

$ws.Range("A512").Value = 44844
$ws.Range("B512").Value = "KA538302"
$ws.Range("C512").Value = "I20"
$ws.Range("D512").Value = "RUNNING REPAIR"
$ws.Range("E512").Value = "WORK DONE DELIVERED"
$ws.Range("F512").Value = 9653
$ws.Range("G512").Value = "P PAY"

$ws.Range("A513").Value = 44845
$ws.Range("B513").Value = "AP09BX8688"
$ws.Range("C513").Value = "RITZ"
$ws.Range("D513").Value = "ABS SENSOR"
$ws.Range("E513").Value = "WORK DONE DELIVERED"
$ws.Range("F513").Value = 1500
$ws.Range("G513").Value = "CREDIT"

$ws.Range("A514").Value = 44845
$ws.Range("B514").Value = "KA53MB5174"
$ws.Range("C514").Value = "XCENT"
$ws.Range("D514").Value = "PMS"
$ws.Range("E514").Value = "WORK DONE DELIVERED"
$ws.Range("F514").Value = 7284
$ws.Range("G514").Value = "CREDIT"

$ws.Range("A515").Value = 44845
$ws.Range("B515").Value = "KA03MQ9315"
$ws.Range("C515").Value = "INNOVA"
$ws.Range("D515").Value = "PMS"
$ws.Range("E515").Value = "WORK DONE DELIVERED"
$ws.Range("F515").Value = 5699
$ws.Range("G515").Value = "G PAY"

$ws.Range("A516").Value = 44846
$ws.Range("B516").Value = "KA03MM9606"
$ws.Range("C516").Value = "H-CITY"
$ws.Range("D516").Value = "PMS"
$ws.Range("E516").Value = "WORK DONE DELIVERED"
$ws.Range("F516").Value = 6938
$ws.Range("G516").Value = "G PAY"

$ws.Range("A517").Value = 44846
$ws.Range("B517").Value = "KA01MV7543"
$ws.Range("C517").Value = "BEAT"
$ws.Range("D517").Value = "GENERAL CHECKUP"
$ws.Range("E517").Value = "WORK DONE DELIVERED"
$ws.Range("F517").Value = 6372
$ws.Range("G517").Value = "G PAY"

$ws.Range("A518").Value = 44846
$ws.Range("B518").Value = "KA01MM2572"
$ws.Range("C518").Value = "INNOVA"
$ws.Range("D518").Value = "RUNNING REPAIR"
$ws.Range("E518").Value = "WORK DONE DELIVERED"
$ws.Range("F518").Value = 9681
$ws.Range("G518").Value = "CREDIT"

$ws.Range("A519").Value = 44846
$ws.Range("B519").Value = "KA04MM681"
$ws.Range("C519").Value = "SPARK"
$ws.Range("D519").Value = "CLUTCH PROBLEM"
$ws.Range("E519").Value = "WORK DONE DELIVERED"
$ws.Range("F519").Value = 15751
$ws.Range("G519").Value = "P PAY"

# --- Update the sheet's used-range dimension / view to reflect the new extent ---
$ws.Range("G519").Select()
